$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "2021" column: copy the formatting from the neighboring 2020 (Q)
# column so the header (row 4) and data (row 5) cells pick up the same
# styles (right-aligned header style, bottom-border data style), then set
# the new values.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 42.9

# Update the sheet's recorded selection to reflect the new last column,
# matching the view state after the column was added.
$ws.Range("R9").Select()
